$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.74
$ws.Range("G2").Value = 2.84
$ws.Range("H2").Value = 2.52
$ws.Range("N2").Value = 5.6
$ws.Range("P2").Value = 2.62
$ws.Range("Q2").Value = 1.6
$ws.Range("R2").Value = 1.62
$ws.Range("S2").Value = 2.52
$ws.Range("Y2").Value = 16
$ws.Range("AA2").Value = 38
$ws.Range("AB2").Value = 17
$ws.Range("AC2").Value = 10
$ws.Range("AE2").Value = 24
$ws.Range("AF2").Value = 22
$ws.Range("AJ2").Value = 42
$ws.Range("AK2").Value = 26
$ws.Range("AL2").Value = 32
$ws.Range("AM2").Value = 60
$ws.Range("AN2").Value = 16
$ws.Range("AO2").Value = 14
$ws.Range("F3").Value = 1.67
$ws.Range("G3").Value = 1.8
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 5.3
$ws.Range("J3").Value = 4.3
$ws.Range("P3").Value = 2.96
$ws.Range("Q3").Value = 1.35
$ws.Range("R3").Value = 1.93
$ws.Range("S3").Value = 1.87
$ws.Range("U3").Value = 2.84
$ws.Range("V3").Value = 1.23
$ws.Range("W3").Value = 2.24
$ws.Range("X3").Value = 40
$ws.Range("Y3").Value = 38
$ws.Range("Z3").Value = 50
$ws.Range("AA3").Value = 110
$ws.Range("AB3").Value = 21
$ws.Range("AC3").Value = 14
$ws.Range("AE3").Value = 46
$ws.Range("AF3").Value = 19.5
$ws.Range("AG3").Value = 12.5
$ws.Range("AH3").Value = 17.5
$ws.Range("AI3").Value = 42
$ws.Range("AJ3").Value = 25
$ws.Range("AK3").Value = 17.5
$ws.Range("AL3").Value = 24
$ws.Range("AN3").Value = 6.4
$ws.Range("AO3").Value = 25
$ws.Range("K4").Value = 4.7
$ws.Range("L4").Value = 1.39
$ws.Range("P4").Value = 1.89
$ws.Range("X4").Value = 15.5
$ws.Range("G5").Value = 3.4
$ws.Range("H5").Value = 2.3
$ws.Range("I5").Value = 3.1
$ws.Range("K5").Value = 5.9
$ws.Range("N5").Value = 1.11
$ws.Range("Q5").Value = 2.12
$ws.Range("S5").Value = 2.12
$ws.Range("T5").Value = 1.03
$ws.Range("U5").Value = 1.03
$ws.Range("V5").Value = 1.47
$ws.Range("W5").Value = 1.41
$ws.Range("G6").Value = 3.85
$ws.Range("L6").Value = 1.34
$ws.Range("W6").Value = 1.36
$ws.Range("G7").Value = 1.37
$ws.Range("H7").Value = 10
$ws.Range("J7").Value = 5.6
$ws.Range("L7").Value = 1.26
$ws.Range("N7").Value = 5.2
$ws.Range("O7").Value = 1.19
$ws.Range("P7").Value = 2.44
$ws.Range("Q7").Value = 1.55
$ws.Range("R7").Value = 1.58
$ws.Range("S7").Value = 2.38
$ws.Range("T7").Value = 1.93
$ws.Range("U7").Value = 1.88
$ws.Range("V7").Value = 1.08
$ws.Range("W7").Value = 3.7
$ws.Range("X7").Value = 32
$ws.Range("Y7").Value = 48
$ws.Range("AB7").Value = 11
$ws.Range("AC7").Value = 16.5
$ws.Range("AD7").Value = 48
$ws.Range("AF7").Value = 10.5
$ws.Range("AG7").Value = 11.5
$ws.Range("AJ7").Value = 12
$ws.Range("AK7").Value = 17
$ws.Range("AL7").Value = 42
$ws.Range("AN7").Value = 5.7
$ws.Range("J8").Value = 3.65
$ws.Range("K8").Value = 3.75
$ws.Range("P8").Value = 2.04
$ws.Range("AO8").Value = 22
$ws.Range("F9").Value = 3.85
$ws.Range("H9").Value = 1.89
$ws.Range("I9").Value = 1.93
$ws.Range("J9").Value = 4.1
$ws.Range("K9").Value = 4.5
$ws.Range("Q9").Value = 1.57
$ws.Range("U9").Value = 2.52
$ws.Range("V9").Value = 2.06
$ws.Range("W9").Value = 1.31
$ws.Range("Y9").Value = 13.5
$ws.Range("AE9").Value = 21
$ws.Range("AJ9").Value = 90
$ws.Range("AM9").Value = 75
$ws.Range("AO9").Value = 8.800000000000001
$ws.Range("G10").Value = 2.56
$ws.Range("I10").Value = 3.5
$ws.Range("L10").Value = 1.36
$ws.Range("M10").Value = 1.08
$ws.Range("P10").Value = 1.83
$ws.Range("S10").Value = 3.6
$ws.Range("T10").Value = 1.78
$ws.Range("X10").Value = 14.5
$ws.Range("Y10").Value = 14
$ws.Range("Z10").Value = 23
$ws.Range("AH10").Value = 18
$ws.Range("AI10").Value = 55
$ws.Range("AN10").Value = 23
$ws.Range("AO10").Value = 46
$ws.Range("G11").Value = 5
$ws.Range("H11").Value = 2
$ws.Range("J11").Value = 1.25
$ws.Range("N11").Value = 1.11
$ws.Range("T11").Value = 1.03
$ws.Range("U11").Value = 1.03
$ws.Range("G12").Value = 3.25
$ws.Range("H12").Value = 2.18
$ws.Range("I12").Value = 2.38
$ws.Range("J12").Value = 3.95
$ws.Range("K12").Value = 4.4
$ws.Range("L12").Value = 1.25
$ws.Range("N12").Value = 5.4
$ws.Range("O12").Value = 1.18
$ws.Range("P12").Value = 2.5
$ws.Range("Q12").Value = 1.54
$ws.Range("R12").Value = 1.6
$ws.Range("S12").Value = 2.38
$ws.Range("V12").Value = 1.75
$ws.Range("AC12").Value = 1000
$ws.Range("AG12").Value = 17
$ws.Range("AH12").Value = 19
$ws.Range("AN12").Value = 23
$ws.Range("F13").Value = 2.76
$ws.Range("H13").Value = 2.78
$ws.Range("I13").Value = 2.88
$ws.Range("J13").Value = 3.4
$ws.Range("O13").Value = 1.34
$ws.Range("F14").Value = 1.49
$ws.Range("Q14").Value = 1.98
$ws.Range("V14").Value = 1.1
$ws.Range("N15").Value = 1.11
$ws.Range("T15").Value = 1.03
$ws.Range("U15").Value = 1.03
$ws.Range("H16").Value = 2.4
$ws.Range("I16").Value = 2.52
$ws.Range("Q16").Value = 3.55
$ws.Range("S16").Value = 8.6
$ws.Range("U16").Value = 1.57
$ws.Range("V16").Value = 1.65
$ws.Range("W16").Value = 1.34
$ws.Range("Z16").Value = 13.5
$ws.Range("AA16").Value = 85
$ws.Range("AD16").Value = 18
$ws.Range("AE16").Value = 150
$ws.Range("AO16").Value = 330
$ws.Range("G17").Value = 2.12
$ws.Range("I17").Value = 5.3
$ws.Range("J17").Value = 3
$ws.Range("K17").Value = 3.4
$ws.Range("N17").Value = 2.48
$ws.Range("O17").Value = 1.56
$ws.Range("T17").Value = 2.22
$ws.Range("U17").Value = 1.7
$ws.Range("W17").Value = 1.89
$ws.Range("AC17").Value = 1000
$ws.Range("H18").Value = 1.85
$ws.Range("I18").Value = 2.12
$ws.Range("K18").Value = 4.1
$ws.Range("O18").Value = 1.53
$ws.Range("T18").Value = 2.22
$ws.Range("V18").Value = 1.89
$ws.Range("W18").Value = 1.21
$ws.Range("Z18").Value = 11.5
$ws.Range("AF18").Value = 46
$ws.Range("AG18").Value = 26
$ws.Range("F19").Value = 1.75
$ws.Range("I19").Value = 5.6
$ws.Range("K19").Value = 4.1
$ws.Range("L19").Value = 1.44
$ws.Range("M19").Value = 1.08
$ws.Range("N19").Value = 3.35
$ws.Range("P19").Value = 1.8
$ws.Range("Q19").Value = 2.04
$ws.Range("R19").Value = 1.3
$ws.Range("U19").Value = 1.9
$ws.Range("V19").Value = 1.21
$ws.Range("U20").Value = 1.9
